# Updated cryptos list on Fri May 31 05:11:04 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell $ws.Range("D2") "68.459.48"
Set-TextCell $ws.Range("E2") "  +0.52%  "

# Row 3 - Ethereum
Set-TextCell $ws.Range("D3") "3.757.52"
Set-TextCell $ws.Range("E3") "  -0.53%  "

# Row 4 - TetherUSD
Set-TextCell $ws.Range("E4") "  -0.01%  "

# Row 5 - BNB
Set-TextCell $ws.Range("D5") "593.12"
Set-TextCell $ws.Range("E5") "  -0.60%  "

# Row 6 - Solana
Set-TextCell $ws.Range("D6") "166.83"
Set-TextCell $ws.Range("E6") "  -1.55%  "

# Row 7 - LidoStakedEther
Set-TextCell $ws.Range("D7") "3.751.80"
Set-TextCell $ws.Range("E7") "  -0.64%  "

# Row 8 - USDC
Set-TextCell $ws.Range("E8") "  +0.00%  "

# Row 9 - XRP
Set-TextCell $ws.Range("E9") "  -1.06%  "

# Row 10 - Dogecoin
Set-TextCell $ws.Range("D10") "0.159"

# Row 11 - Toncoin
Set-TextCell $ws.Range("D11") "6.40"
Set-TextCell $ws.Range("E11") "  -1.52%  "

# Row 12 - Cardano
Set-TextCell $ws.Range("D12") "0.448"
Set-TextCell $ws.Range("E12") "  -0.82%  "

# Row 13 - ShibaInu
Set-TextCell $ws.Range("D13") "0.0000259"
Set-TextCell $ws.Range("E13") "  -6.53%  "

# Row 14 - Avalanche
Set-TextCell $ws.Range("D14") "36.05"
Set-TextCell $ws.Range("E14") "  -1.62%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextCell $ws.Range("D15") "4.390.47"
Set-TextCell $ws.Range("E15") "  -0.52%  "

# Row 16 - WrappedEther
Set-TextCell $ws.Range("D16") "3.757.07"
Set-TextCell $ws.Range("E16") "  -0.40%  "

# Row 17 - WrappedBTC
Set-TextCell $ws.Range("D17") "68.453.72"
Set-TextCell $ws.Range("E17") "  +0.52%  "

# Row 18 - Chainlink
Set-TextCell $ws.Range("D18") "17.96"
Set-TextCell $ws.Range("E18") "  -3.93%  "

# Row 19 - TRON
Set-TextCell $ws.Range("E19") "  +0.86%  "

# Row 20 - Polkadot
Set-TextCell $ws.Range("D20") "6.98"
Set-TextCell $ws.Range("E20") "  -2.53%  "

# Row 21 - Uniswap
Set-TextCell $ws.Range("D21") "10.69"
Set-TextCell $ws.Range("E21") "  +1.39%  "

# Row 22 - BitcoinCash
Set-TextCell $ws.Range("D22") "465.41"
Set-TextCell $ws.Range("E22") "  -0.31%  "

# Row 23 - Polygon
Set-TextCell $ws.Range("D23") "0.696"
Set-TextCell $ws.Range("E23") "  -2.87%  "

# Row 24 - PEPE -> Litecoin (swap)
Set-TextCell $ws.Range("B24") "Litecoin"
Set-TextCell $ws.Range("C24") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell $ws.Range("D24") "84.00"
Set-TextCell $ws.Range("E24") "  +0.16%  "

# Row 25 - Litecoin -> PEPE (swap)
Set-TextCell $ws.Range("B25") "PEPE"
Set-TextCell $ws.Range("C25") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws.Range("D25") "0.0000147"
Set-TextCell $ws.Range("E25") "  -1.05%  "

# Row 26 - Fetch.AI
Set-TextCell $ws.Range("D26") "2.18"
Set-TextCell $ws.Range("E26") "  -2.64%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextCell $ws.Range("D27") "11.90"
Set-TextCell $ws.Range("E27") "  -1.81%  "

# Row 28 - RenderToken
Set-TextCell $ws.Range("D28") "10.05"
Set-TextCell $ws.Range("E28") "  -2.95%  "

# Row 29 - Dai
Set-TextCell $ws.Range("E29") "  -0.07%  "

# Row 30 - WrappedeETH
Set-TextCell $ws.Range("D30") "3.908.19"
Set-TextCell $ws.Range("E30") "  -0.54%  "

# Row 31 - PancakeSwap
Set-TextCell $ws.Range("E31") "  -4.55%  "

# Row 32 - NEARProtocol
Set-TextCell $ws.Range("D32") "7.32"
Set-TextCell $ws.Range("E32") "  -3.23%  "

# Row 33 - EthereumClassic
Set-TextCell $ws.Range("D33") "29.93"
Set-TextCell $ws.Range("E33") "  -1.75%  "

# Row 34 - ImmutableX
Set-TextCell $ws.Range("D34") "2.17"
Set-TextCell $ws.Range("E34") "  -2.31%  "

# Row 35 - Aptos
Set-TextCell $ws.Range("D35") "9.17"
Set-TextCell $ws.Range("E35") "  -0.34%  "

# Row 36 - Binance-PegBSC-USD
Set-TextCell $ws.Range("D36") "0.999"

# Row 37 - RenzoRestakedETH
Set-TextCell $ws.Range("D37") "3.716.25"
Set-TextCell $ws.Range("E37") "  -0.53%  "

# Row 38 - Hedera
Set-TextCell $ws.Range("D38") "0.100"
Set-TextCell $ws.Range("E38") "  -3.08%  "

# Row 39 - dogwifhat
Set-TextCell $ws.Range("D39") "3.38"
Set-TextCell $ws.Range("E39") "  -9.98%  "

# Row 40 - Mantle
Set-TextCell $ws.Range("D40") "0.999"
Set-TextCell $ws.Range("E40") "  -0.19%  "

# Row 41 - Kaspa
Set-TextCell $ws.Range("E41") "  -1.25%  "

# Row 42 - Filecoin
Set-TextCell $ws.Range("D42") "5.78"
Set-TextCell $ws.Range("E42") "  -1.38%  "

# Row 43 - FirstDigitalUSD
Set-TextCell $ws.Range("E43") "  -0.03%  "

# Row 44 - USDe
Set-TextCell $ws.Range("E44") "  +0.00%  "

# Row 45 - TheGraph
Set-TextCell $ws.Range("D45") "0.302"
Set-TextCell $ws.Range("E45") "  -2.99%  "

# Row 46 - Arweave
Set-TextCell $ws.Range("D46") "43.81"
Set-TextCell $ws.Range("E46") "  +9.41%  "

# Row 47 - OKB
Set-TextCell $ws.Range("D47") "46.64"
Set-TextCell $ws.Range("E47") "  +2.49%  "

# Row 48 - Stacks
Set-TextCell $ws.Range("E48") "  -1.81%  "

# Row 49 - Cosmos
Set-TextCell $ws.Range("D49") "8.47"
Set-TextCell $ws.Range("E49") "  -2.07%  "

# Row 50 - Monero
Set-TextCell $ws.Range("D50") "145.60"
Set-TextCell $ws.Range("E50") "  +1.45%  "

# Row 51 - Bittensor
Set-TextCell $ws.Range("D51") "390.01"
Set-TextCell $ws.Range("E51") "  -3.60%  "
